$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.729584333333333
$ws.Range("H2").Value = 5.188753
$ws.Range("I2").Value = 0.2476387648475193
$ws.Range("J2").Value = 0.2476387648475193
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 12.87733
$ws.Range("N2").Value = 38.63199
$ws.Range("O2").Value = 0.1584922499374361
$ws.Range("P2").Value = 0.1584922499374361
$ws.Range("Q2").Value = 22.27242822316333
$ws.Range("R2").Value = 200.45185400847
$ws.Range("S2").Value = 0.039248825012411
$ws.Range("T2").Value = 0.039248825012411
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.729584333333333
$ws.Range("H3").Value = 5.188753
$ws.Range("I3").Value = 0.2476387648475193
$ws.Range("J3").Value = 0.2476387648475193
$ws.Range("O3").Value = 0.4359831802722915
$ws.Range("P3").Value = 0.4359831802722916
$ws.Range("Q3").Value = 61.26737485873422
$ws.Range("R3").Value = 551.4063737286081
$ws.Range("S3").Value = 0.1079663362569236
$ws.Range("T3").Value = 0.1079663362569236
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.729584333333333
$ws.Range("H4").Value = 5.188753
$ws.Range("I4").Value = 0.2476387648475193
$ws.Range("J4").Value = 0.2476387648475193
$ws.Range("M4").Value = 30.51453966666667
$ws.Range("N4").Value = 91.54361900000001
$ws.Range("O4").Value = 0.3755683862706898
$ws.Range("P4").Value = 0.3755683862706898
$ws.Range("Q4").Value = 52.77746974634523
$ws.Range("R4").Value = 474.997227717107
$ws.Range("S4").Value = 0.09300529129184963
$ws.Range("T4").Value = 0.09300529129184965
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.729584333333333
$ws.Range("H5").Value = 5.188753
$ws.Range("I5").Value = 0.2476387648475193
$ws.Range("J5").Value = 0.2476387648475193
$ws.Range("M5").Value = 2.433908666666667
$ws.Range("N5").Value = 7.301726
$ws.Range("O5").Value = 0.0299561835195825
$ws.Range("P5").Value = 0.0299561835195825
$ws.Range("Q5").Value = 4.20965029863089
$ws.Range("R5").Value = 37.886852687678
$ws.Range("S5").Value = 0.007418312286335022
$ws.Range("T5").Value = 0.007418312286335022
$ws.Range("I6").Value = 0.2307941364328804
$ws.Range("J6").Value = 0.2307941364328804
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 12.87733
$ws.Range("N6").Value = 38.63199
$ws.Range("O6").Value = 0.1584922499374361
$ws.Range("P6").Value = 0.1584922499374361
$ws.Range("Q6").Value = 20.75743610332334
$ws.Range("R6").Value = 186.81692492991
$ws.Range("S6").Value = 0.03657908195561481
$ws.Range("T6").Value = 0.03657908195561481
$ws.Range("I7").Value = 0.2307941364328804
$ws.Range("J7").Value = 0.2307941364328804
$ws.Range("O7").Value = 0.4359831802722915
$ws.Range("P7").Value = 0.4359831802722916
$ws.Range("S7").Value = 0.1006223615902043
$ws.Range("T7").Value = 0.1006223615902044
$ws.Range("I8").Value = 0.2307941364328804
$ws.Range("J8").Value = 0.2307941364328804
$ws.Range("M8").Value = 30.51453966666667
$ws.Range("N8").Value = 91.54361900000001
$ws.Range("O8").Value = 0.3755683862706898
$ws.Range("P8").Value = 0.3755683862706898
$ws.Range("Q8").Value = 49.18749518364123
$ws.Range("R8").Value = 442.6874566527711
$ws.Range("S8").Value = 0.0866789813808343
$ws.Range("T8").Value = 0.08667898138083431
$ws.Range("I9").Value = 0.2307941364328804
$ws.Range("J9").Value = 0.2307941364328804
$ws.Range("M9").Value = 2.433908666666667
$ws.Range("N9").Value = 7.301726
$ws.Range("O9").Value = 0.0299561835195825
$ws.Range("P9").Value = 0.0299561835195825
$ws.Range("Q9").Value = 3.92330581181489
$ws.Range("R9").Value = 35.309752306334
$ws.Range("S9").Value = 0.006913711506226925
$ws.Range("T9").Value = 0.006913711506226925
$ws.Range("G10").Value = 2.743651333333334
$ws.Range("H10").Value = 8.230954000000001
$ws.Range("I10").Value = 0.3928310486309039
$ws.Range("J10").Value = 0.3928310486309038
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 12.87733
$ws.Range("N10").Value = 38.63199
$ws.Range("O10").Value = 0.1584922499374361
$ws.Range("P10").Value = 0.1584922499374361
$ws.Range("Q10").Value = 35.33090362427334
$ws.Range("R10").Value = 317.9781326184601
$ws.Range("S10").Value = 0.06226067674279435
$ws.Range("T10").Value = 0.06226067674279434
$ws.Range("G11").Value = 2.743651333333334
$ws.Range("H11").Value = 8.230954000000001
$ws.Range("I11").Value = 0.3928310486309039
$ws.Range("J11").Value = 0.3928310486309038
$ws.Range("O11").Value = 0.4359831802722915
$ws.Range("P11").Value = 0.4359831802722916
$ws.Range("Q11").Value = 97.1888513797049
$ws.Range("R11").Value = 874.699662417344
$ws.Range("S11").Value = 0.1712677298918007
$ws.Range("T11").Value = 0.1712677298918007
$ws.Range("G12").Value = 2.743651333333334
$ws.Range("H12").Value = 8.230954000000001
$ws.Range("I12").Value = 0.3928310486309039
$ws.Range("J12").Value = 0.3928310486309038
$ws.Range("M12").Value = 30.51453966666667
$ws.Range("N12").Value = 91.54361900000001
$ws.Range("O12").Value = 0.3755683862706898
$ws.Range("P12").Value = 0.3755683862706898
$ws.Range("Q12").Value = 83.7212574425029
$ws.Range("R12").Value = 753.4913169825261
$ws.Range("S12").Value = 0.1475349230113314
$ws.Range("T12").Value = 0.1475349230113314
$ws.Range("G13").Value = 2.743651333333334
$ws.Range("H13").Value = 8.230954000000001
$ws.Range("I13").Value = 0.3928310486309039
$ws.Range("J13").Value = 0.3928310486309038
$ws.Range("M13").Value = 2.433908666666667
$ws.Range("N13").Value = 7.301726
$ws.Range("O13").Value = 0.0299561835195825
$ws.Range("P13").Value = 0.0299561835195825
$ws.Range("Q13").Value = 6.677796758511557
$ws.Range("R13").Value = 60.10017082660401
$ws.Range("S13").Value = 0.01176771898497739
$ws.Range("T13").Value = 0.01176771898497739
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.8991316666666668
$ws.Range("H14").Value = 2.697395
$ws.Range("I14").Value = 0.1287360500886965
$ws.Range("J14").Value = 0.1287360500886965
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 12.87733
$ws.Range("N14").Value = 38.63199
$ws.Range("O14").Value = 0.1584922499374361
$ws.Range("P14").Value = 0.1584922499374361
$ws.Range("Q14").Value = 11.57841518511667
$ws.Range("R14").Value = 104.20573666605
$ws.Range("S14").Value = 0.02040366622661598
$ws.Range("T14").Value = 0.02040366622661598
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.8991316666666668
$ws.Range("H15").Value = 2.697395
$ws.Range("I15").Value = 0.1287360500886965
$ws.Range("J15").Value = 0.1287360500886965
$ws.Range("O15").Value = 0.4359831802722915
$ws.Range("P15").Value = 0.4359831802722916
$ws.Range("Q15").Value = 31.85010167319111
$ws.Range("R15").Value = 286.65091505872
$ws.Range("S15").Value = 0.05612675253336293
$ws.Range("T15").Value = 0.05612675253336293
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.8991316666666668
$ws.Range("H16").Value = 2.697395
$ws.Range("I16").Value = 0.1287360500886965
$ws.Range("J16").Value = 0.1287360500886965
$ws.Range("M16").Value = 30.51453966666667
$ws.Range("N16").Value = 91.54361900000001
$ws.Range("O16").Value = 0.3755683862706898
$ws.Range("P16").Value = 0.3755683862706898
$ws.Range("Q16").Value = 27.43658890805612
$ws.Range("R16").Value = 246.929300172505
$ws.Range("S16").Value = 0.04834919058667444
$ws.Range("T16").Value = 0.04834919058667445
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.8991316666666668
$ws.Range("H17").Value = 2.697395
$ws.Range("I17").Value = 0.1287360500886965
$ws.Range("J17").Value = 0.1287360500886965
$ws.Range("M17").Value = 2.433908666666667
$ws.Range("N17").Value = 7.301726
$ws.Range("O17").Value = 0.0299561835195825
$ws.Range("P17").Value = 0.0299561835195825
$ws.Range("Q17").Value = 2.188404355974445
$ws.Range("R17").Value = 19.69563920377
$ws.Range("S17").Value = 0.003856440742043157
$ws.Range("T17").Value = 0.003856440742043157
